$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 1703
$ws1.Range("F12").Value = 1405
$ws1.Range("F16").Value = 12675
$ws1.Range("F17").Value = 12697
$ws1.Range("F18").Value = 942
$ws1.Range("F23").Value = 517
$ws1.Range("F24").Value = 1988
$ws1.Range("F27").Value = 234

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 68

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 1703
$ws4.Range("F17").Value = 1405
$ws4.Range("F22").Value = 12675
$ws4.Range("F23").Value = 12698
$ws4.Range("F24").Value = 942
$ws4.Range("F29").Value = 517
$ws4.Range("F32").Value = 1988
$ws4.Range("F37").Value = 234
$ws4.Range("F39").Value = 68
